# Weekly update: insert a new price record as the new first row of the
# "Vega Modelo de Temuco" / Coco weekly series (row 69), pushing the
# previously-existing rows 69-77 down to 70-78.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 69; existing rows 69:77 shift to 70:78.
$ws.Rows.Item(69).Insert()

# Populate the newly inserted row 69 with this week's data.
$ws.Cells.Item(69, 1).Value  = 10
$ws.Cells.Item(69, 2).Value  = 'Vega Modelo de Temuco'
$ws.Cells.Item(69, 3).Value  = 'La Araucanía'
$ws.Cells.Item(69, 4).Value  = 44746
$ws.Cells.Item(69, 5).Value  = 9
$ws.Cells.Item(69, 6).Value  = 'Fruta'
$ws.Cells.Item(69, 7).Value  = 100108
$ws.Cells.Item(69, 8).Value  = 'Tropicales y subtropicales'
$ws.Cells.Item(69, 9).Value  = 100108007
$ws.Cells.Item(69, 10).Value = 'Coco'
$ws.Cells.Item(69, 11).Value = 'Sin especificar'
$ws.Cells.Item(69, 12).Value = 'Primera'
$ws.Cells.Item(69, 13).Value = 70
$ws.Cells.Item(69, 14).Value = 28000
$ws.Cells.Item(69, 15).Value = 30000
$ws.Cells.Item(69, 16).Value = 29143
$ws.Cells.Item(69, 17).Value = '$/malla 20 unidades'
$ws.Cells.Item(69, 18).Value = 'Perú'
$ws.Cells.Item(69, 19).Value = 1457
$ws.Cells.Item(69, 20).Value = 20
